$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 29-44 (old extra data rows), shrinking dimension to A1:C28
$ws.Rows("29:44").Delete() | Out-Null

# Update data for rows 2-28 with new values (location changed from Astro to Dabaca)
$ws.Range("A2").Value2 = 45807
$ws.Range("B2").Value2 = 12.187
$ws.Range("C2").Value2 = 0
$ws.Range("A3").Value2 = 45807.01041666666
$ws.Range("B3").Value2 = 16.47
$ws.Range("C3").Value2 = 0
$ws.Range("A4").Value2 = 45807.02083333334
$ws.Range("B4").Value2 = 14.575
$ws.Range("C4").Value2 = 0
$ws.Range("A5").Value2 = 45807.03125
$ws.Range("B5").Value2 = 3.637
$ws.Range("C5").Value2 = 1.129
$ws.Range("A6").Value2 = 45807.04166666666
$ws.Range("B6").Value2 = 5.038
$ws.Range("C6").Value2 = 1.118
$ws.Range("A7").Value2 = 45807.05208333334
$ws.Range("B7").Value2 = 1.29
$ws.Range("C7").Value2 = 2.53
$ws.Range("A8").Value2 = 45807.0625
$ws.Range("B8").Value2 = 0
$ws.Range("C8").Value2 = 6.198
$ws.Range("A9").Value2 = 45807.07291666666
$ws.Range("B9").Value2 = 0
$ws.Range("C9").Value2 = 14.092
$ws.Range("A10").Value2 = 45807.08333333334
$ws.Range("B10").Value2 = 0
$ws.Range("C10").Value2 = 14.61
$ws.Range("A11").Value2 = 45807.09375
$ws.Range("B11").Value2 = 0
$ws.Range("C11").Value2 = 11.589
$ws.Range("A12").Value2 = 45807.10416666666
$ws.Range("B12").Value2 = 0
$ws.Range("C12").Value2 = 4.1
$ws.Range("A13").Value2 = 45807.11458333334
$ws.Range("B13").Value2 = 0
$ws.Range("C13").Value2 = 11.145
$ws.Range("A14").Value2 = 45807.125
$ws.Range("B14").Value2 = 0
$ws.Range("C14").Value2 = 9.952999999999999
$ws.Range("A15").Value2 = 45807.13541666666
$ws.Range("B15").Value2 = 0
$ws.Range("C15").Value2 = 1.711
$ws.Range("A16").Value2 = 45807.14583333334
$ws.Range("B16").Value2 = 0
$ws.Range("C16").Value2 = 1.234
$ws.Range("A17").Value2 = 45807.15625
$ws.Range("B17").Value2 = 0
$ws.Range("C17").Value2 = 1.965
$ws.Range("A18").Value2 = 45807.16666666666
$ws.Range("B18").Value2 = 0.017
$ws.Range("C18").Value2 = 2.121
$ws.Range("A19").Value2 = 45807.17708333334
$ws.Range("B19").Value2 = 0
$ws.Range("C19").Value2 = 2.379
$ws.Range("A20").Value2 = 45807.1875
$ws.Range("B20").Value2 = 0.02
$ws.Range("C20").Value2 = 5.745
$ws.Range("A21").Value2 = 45807.19791666666
$ws.Range("B21").Value2 = 0
$ws.Range("C21").Value2 = 4.214
$ws.Range("A22").Value2 = 45807.20833333334
$ws.Range("B22").Value2 = 6.28
$ws.Range("C22").Value2 = 0.065
$ws.Range("A23").Value2 = 45807.21875
$ws.Range("B23").Value2 = 8.981
$ws.Range("C23").Value2 = 0
$ws.Range("A24").Value2 = 45807.22916666666
$ws.Range("B24").Value2 = 1.886
$ws.Range("C24").Value2 = 0.052
$ws.Range("A25").Value2 = 45807.23958333334
$ws.Range("B25").Value2 = 13.915
$ws.Range("C25").Value2 = 0.08500000000000001
$ws.Range("A26").Value2 = 45807.25
$ws.Range("B26").Value2 = 15.104
$ws.Range("C26").Value2 = 0.002
$ws.Range("A27").Value2 = 45807.26041666666
$ws.Range("B27").Value2 = 18.041
$ws.Range("C27").Value2 = 0
$ws.Range("A28").Value2 = 45807.27083333334
$ws.Range("B28").Value2 = 46.775
$ws.Range("C28").Value2 = 0
